$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data rows 2-11: employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary
$data = @(
    @(80366, "Sr. Luiz Miguel Carvalho", "Financeiro", "Problemas pessoais", 7, 45078, 4320.45),
    @(16059, "Clara Costa", "Atendimento ao Cliente", "Outros", 1, 45093, 2203.16),
    @(74319, "Dr. Gael Moura", "Recursos Humanos", "Doenca", 8, 45104, 5662.22),
    @(47304, "Maria Alice Macedo", "TI", "Problemas pessoais", 4, 45083, 8922.530000000001),
    @(29914, "Sr. Igor da Rocha", "Recursos Humanos", "Viagem de negocios", 6, 45078, 3713.32),
    @(96546, "Helena Macedo", "Atendimento ao Cliente", "Viagem de negocios", 2, 45086, 4262.92),
    @(76429, "Kaique da Luz", "Operacoes", "Problemas pessoais", 5, 45090, 5975.62),
    @(13015, "Vinícius Camargo", "Vendas", "Consulta medica", 1, 45092, 6848.67),
    @(15388, "João Lucas Andrade", "Recursos Humanos", "Problemas pessoais", 4, 45096, 3889.13),
    @(80453, "João Lucas Cardoso", "Juridico", "Doenca", 8, 45083, 8909.290000000001)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rec[6]
    $rowIndex++
}
